$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Ref, $Val)
    $c = $ws.Range($Ref)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

Set-TextCell "D2" "27.721.41"
Set-TextCell "D3" "1.851.34"
Set-TextCell "E3" "  +0.25%  "
Set-TextCell "E4" "  +0.63%  "
Set-TextCell "D5" "322.38"
Set-TextCell "E5" "  +0.61%  "
Set-TextCell "D6" "1.030"
Set-TextCell "E6" "  +0.47%  "
Set-TextCell "D7" "0.4393"
Set-TextCell "E7" "  +0.59%  "
Set-TextCell "D8" "0.3807"
Set-TextCell "E8" "  +1.79%  "
Set-TextCell "D9" "0.07409"
Set-TextCell "D10" "0.8856"
Set-TextCell "E10" "  +1.24%  "
Set-TextCell "D11" "21.56"
Set-TextCell "E11" "  +0.49%  "
Set-TextCell "D12" "1.859.83"
Set-TextCell "E12" "  +0.30%  "
Set-TextCell "D13" "5.509"
Set-TextCell "E13" "  +0.29%  "
Set-TextCell "D14" "6.708"
Set-TextCell "E14" "  +0.51%  "
Set-TextCell "D15" "0.07179"
Set-TextCell "E15" "  +0.41%  "
Set-TextCell "D16" "85.31"
Set-TextCell "E16" "  +3.07%  "
Set-TextCell "E17" "  +0.84%  "
Set-TextCell "D18" "0.000009052"
Set-TextCell "E18" "  +0.40%  "
Set-TextCell "D19" "1.031"
Set-TextCell "E19" "  +0.52%  "
Set-TextCell "D20" "15.51"
Set-TextCell "E20" "  +0.69%  "
Set-TextCell "D21" "27.759.38"
Set-TextCell "E21" "  +0.54%  "
Set-TextCell "D22" "5.281"
Set-TextCell "E22" "  +0.72%  "
Set-TextCell "D23" "11.25"
Set-TextCell "E23" "  +0.41%  "
Set-TextCell "D24" "2.086.60"
Set-TextCell "E24" "  +0.76%  "
Set-TextCell "D25" "2.056"
Set-TextCell "E25" "  +6.66%  "
Set-TextCell "D26" "159.14"
Set-TextCell "E26" "  +1.38%  "
Set-TextCell "D27" "18.72"
Set-TextCell "E27" "  +0.08%  "
Set-TextCell "D28" "1.997"
Set-TextCell "E28" "  +2.30%  "
Set-TextCell "D29" "5.351"
Set-TextCell "E29" "  +1.27%  "
Set-TextCell "D30" "117.94"
Set-TextCell "E30" "  +1.51%  "
Set-TextCell "D31" "0.09078"
Set-TextCell "E31" "  +0.16%  "
Set-TextCell "D32" "1.211"
Set-TextCell "E32" "  +0.28%  "
Set-TextCell "D33" "0.7717"
Set-TextCell "E33" "  +0.65%  "
Set-TextCell "D34" "3.008"
Set-TextCell "E34" "  +4.77%  "
Set-TextCell "D35" "4.580"
Set-TextCell "E35" "  +1.65%  "
Set-TextCell "D36" "1.032"
Set-TextCell "E36" "  +0.53%  "
Set-TextCell "D37" "1.151"
Set-TextCell "E37" "  +0.38%  "
Set-TextCell "E38" "  +0.11%  "
Set-TextCell "D39" "0.05283"
Set-TextCell "E39" "  +0.39%  "
Set-TextCell "D40" "2.858"
Set-TextCell "E40" "  +2.12%  "
Set-TextCell "D41" "0.5192"
Set-TextCell "E41" "  +0.57%  "
Set-TextCell "D42" "6.886"
Set-TextCell "E42" "  +2.78%  "
Set-TextCell "D43" "0.1670"
Set-TextCell "E43" "  -0.08%  "
Set-TextCell "D44" "8.737"
Set-TextCell "E44" "  +1.79%  "
Set-TextCell "D45" "110.53"
Set-TextCell "E45" "  +1.50%  "
Set-TextCell "D46" "10.71"
Set-TextCell "E46" "  +1.52%  "
Set-TextCell "D47" "1.032"
Set-TextCell "E47" "  +0.53%  "
Set-TextCell "D48" "0.06573"
Set-TextCell "E48" "  +2.91%  "
Set-TextCell "D49" "1.707"
Set-TextCell "E49" "  -0.27%  "
Set-TextCell "D50" "0.4713"
Set-TextCell "E50" "  +1.33%  "
Set-TextCell "D51" "1.894"
Set-TextCell "E51" "  +0.13%  "

